$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Munka1")

# --- Parts table (rows 6-20): row 13 "Váltókábel" -> "Váltó", price/labor updated ---
$ws.Range("F13").Value = "Váltó"
$ws.Range("G13").Value = 20
$ws.Range("H13").Value = 20
$ws.Range("J13").Clear()
$ws.Range("K13").Clear()

# Row 15 "Pedál" labor cost change 15 -> 5
$ws.Range("H15").Value = 5

# --- Repair table header row (27): two new "debug" notes added ---
$ws.Range("I27").Value = "ez nem kell"
$ws.Range("J27").Value = "ez nem kell"
$ws.Range("M27").Value = "Ez sem kell a táblába"

# --- Fix duplicate REPAIRID 102 and renumber everything below it by +1 ---
$ws.Range("F31").Value = 103
$ws.Range("F32").Value = 104
$ws.Range("F33").Value = 105
$ws.Range("F34").Value = 106
$ws.Range("F35").Value = 107
$ws.Range("F36").Value = 108
$ws.Range("F37").Value = 109
$ws.Range("F38").Value = 110
$ws.Range("F39").Value = 111

# Row 32 (REPAIRID 104): rename repair + add part-change columns
$ws.Range("G32").Value = "Fékbetét Csere"
$ws.Range("J32").Value = 16
$ws.Range("K32").Value = 31
$ws.Range("M32").Value = "'-2 fekbetet"

# Row 37 (REPAIRID 109, Pedál Cseréje): add part-change columns
$ws.Range("J37").Value = 10
$ws.Range("K37").Value = 25
$ws.Range("M37").Value = "'-2 pedal"

# Row 39 (REPAIRID 111, Nyereg Cseréje): update labor cost + add part-change columns
$ws.Range("I39").Value = 10
$ws.Range("J39").Value = 10
$ws.Range("K39").Value = 20
$ws.Range("M39").Value = "'-1 nyereg"

# New row 40: additional repair "Gumi Csere" (REPAIRID 112)
$ws.Range("F40").Value = 112
$ws.Range("G40").Value = "Gumi Csere"
$ws.Range("H40").Value = "Első hátsó külső gumi csere"
$ws.Range("I40").Value = 10
$ws.Range("J40").Value = 20
$ws.Range("K40").Value = 30
$ws.Range("M40").Value = "'-2 kulso gumi"

# Restore previous selection/cursor position on the sheet
$ws.Range("F13").Select()
